$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header cells: "_old" -> "_FV2404", "_new" -> "_FV2410" ---
$ws.Range("A1").Value = "Segmentname_FV2404"
$ws.Range("B1").Value = "Segmentgruppe_FV2404"
$ws.Range("C1").Value = "Segment_FV2404"
$ws.Range("D1").Value = "Datenelement_FV2404"
$ws.Range("E1").Value = "Segment ID_FV2404"
$ws.Range("F1").Value = "Code_FV2404"
$ws.Range("G1").Value = "Qualifier_FV2404"
$ws.Range("H1").Value = "Beschreibung_FV2404"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2404"
$ws.Range("J1").Value = "Bedingung_FV2404"
$ws.Range("L1").Value = "Segmentname_FV2410"
$ws.Range("M1").Value = "Segmentgruppe_FV2410"
$ws.Range("N1").Value = "Segment_FV2410"
$ws.Range("O1").Value = "Datenelement_FV2410"
$ws.Range("P1").Value = "Segment ID_FV2410"
$ws.Range("Q1").Value = "Code_FV2410"
$ws.Range("R1").Value = "Qualifier_FV2410"
$ws.Range("S1").Value = "Beschreibung_FV2410"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2410"
$ws.Range("U1").Value = "Bedingung_FV2410"

# --- Turn the used range into a native Excel Table (ListObject) ---
# Stash the existing header formatting off to the side first: ListObjects.Add
# bakes an existing custom header style into a dxf/table style, which the
# source workbook does not have (headers keep plain cell style "s=1" and the
# table carries no explicit style/dxf).
$ws.Range("A1:U1").Copy()
$ws.Range("W1:AQ1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A1:U1").ClearFormats()

$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U64"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# restore the original header formatting
$ws.Range("W1:AQ1").Copy()
$ws.Range("A1:U1").PasteSpecial(-4122)
$ws.Range("W1:AQ1").Clear()

# --- Freeze the header row ---
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select() | Out-Null
